$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Record the path to the document's associated image (new "image" column
# entry for this row) - this is the actual content edit described by the
# commit ("Added paths to images in documents").
$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Insurance/insurance.jpg"

# Mirror the author's resulting selection/view state (cell P2 becomes the
# active cell/selection after the edit).
[void]$ws.Range("P2").Select()
